$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.096.22'
$ws.Range("E2").Value = '  +3.55%  '
$ws.Range("D3").Value = '3.059.80'
$ws.Range("E3").Value = '  +6.36%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '513.78'
$ws.Range("E5").Value = '  +5.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.25'
$ws.Range("E6").Value = '  +7.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.433'
$ws.Range("E8").Value = '  +4.52%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("E10").Value = '  +5.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.371'
$ws.Range("E11").Value = '  +7.94%  '
$ws.Range("D12").Value = '3.576.59'
$ws.Range("E12").Value = '  +6.00%  '
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.21'
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000164'
$ws.Range("E15").Value = '  +5.70%  '
$ws.Range("D16").Value = '57.133.63'
$ws.Range("E16").Value = '  +3.63%  '
$ws.Range("D17").Value = '3.054.89'
$ws.Range("E17").Value = '  +6.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.92'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.05'
$ws.Range("E19").Value = '  +6.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  +7.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '335.03'
$ws.Range("E21").Value = '  +8.27%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.503'
$ws.Range("E23").Value = '  +6.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.14'
$ws.Range("E24").Value = '  +5.67%  '
$ws.Range("E25").Value = '  +6.23%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '0.0₃0944'
$ws.Range("E27").Value = '  +14.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.42'
$ws.Range("E28").Value = '  +2.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.96'
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("E30").Value = '  +5.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.71'
$ws.Range("E31").Value = '  +6.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.18'
$ws.Range("E32").Value = '  +7.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.45'
$ws.Range("E33").Value = '  +4.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("E34").Value = '  +4.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.84'
$ws.Range("E35").Value = '  +6.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.43'
$ws.Range("E36").Value = '  +9.53%  '
$ws.Range("E37").Value = '  +5.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0668'
$ws.Range("E38").Value = '  +4.07%  '
$ws.Range("D39").Value = '3.091.17'
$ws.Range("E39").Value = '  +6.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.92'
$ws.Range("E40").Value = '  +2.83%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.667'
$ws.Range("E41").Value = '  +6.68%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.82'
$ws.Range("E43").Value = '  +6.14%  '
$ws.Range("D44").Value = '2.231.12'
$ws.Range("E44").Value = '  +7.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0251'
$ws.Range("E45").Value = '  +10.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.37'
$ws.Range("E46").Value = '  +5.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.937'
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.85'
$ws.Range("E48").Value = '  +8.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.83'
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0862'
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.180'
$ws.Range("E51").Value = '  +6.21%  '
